# Applies the BOM update for Echo Cinematic:
#  - J11-J15 (row 20): LCSC part # C59983 -> C3975131
#  - R19, R28 (row 34): designator "R19, R28" -> "R19, R28, R51, R52", QTY 2 -> 4
#  - Remove the U1 / ES_DAISY_SEED2_DFM row entirely (row 49), shifting U2.. rows up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update designator + quantity for the 0R resistor row (now includes R51, R52)
$ws.Range("A34").Value = "R19, R28, R51, R52"
$ws.Range("B34").Value = 4

# 2) Update LCSC part number for J11-J15 connector row
$ws.Range("E20").Value = "C3975131"

# 3) Delete the U1 (ES_DAISY_SEED2_DFM) row entirely, shifting rows below up
$ws.Rows(49).Delete()

# Update selection to mirror the author's final cursor position after the row delete
$ws.Range("A49:XFD49").Select()
